# Apply the "Key to Variables" acronym-key update:
# Insert a new row for "IHDbT" / "Industrial Heat Demand by Temperature"
# right before the existing "MHV" row (row 145), pushing everything below
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a new blank row at position 145 (shifts rows 145+ down to 146+)
$ws.Rows.Item(145).Insert()

# Copy the formatting of the row that now sits just below the new blank row
# (row 146, the old row 145 = "MHV" entry) onto the new row so the new
# entry visually matches its neighbors (fill/wrap/font styles, etc).
$ws.Range("A146:G146").Copy()
$ws.Range("A145").PasteSpecial(-4122)  # xlPasteFormats

# The old row 145 only had cells in columns A, B, C and F - clear the
# extra copied cells (D, E, G) so the new row matches that same shape.
$ws.Range("D145:E145").Clear()
$ws.Range("G145").Clear()

# Populate the new row's contents
$ws.Range("A145").Value = "indst"
$ws.Range("B145").Value = "IHDbT"
$ws.Range("C145").Value = "Industrial Heat Demand by Temperature"
$ws.Range("F145").Value = "low"

# Reflect the author's final selection/active cell in the saved view
$ws.Range("A145").Select()
